$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.141.14"
$ws.Range("E2").Value = "  -0.92%  "

$ws.Range("D3").Value = "1.671.10"
$ws.Range("E3").Value = "  -1.31%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.82%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.75"
$ws.Range("E5").Value = "  -3.85%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5256"
$ws.Range("E6").Value = "  -4.35%  "

$ws.Range("E7").Value = "  -0.77%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2651"
$ws.Range("E8").Value = "  -3.22%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06272"
$ws.Range("E9").Value = "  -2.90%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.17"
$ws.Range("E10").Value = "  -3.71%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07515"

$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.435"
$ws.Range("E12").Value = "  -2.27%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.643.12"
$ws.Range("E13").Value = "  -3.28%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5617"
$ws.Range("E14").Value = "  -3.64%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008005"
$ws.Range("E15").Value = "  -4.11%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.43"
$ws.Range("E16").Value = "  +1.50%  "

$ws.Range("D17").Value = "26.182.49"
$ws.Range("E17").Value = "  -0.93%  "

$ws.Range("E18").Value = "  -0.80%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.791"
$ws.Range("E19").Value = "  -2.97%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "187.57"
$ws.Range("E20").Value = "  -2.35%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.37"
$ws.Range("E21").Value = "  -5.45%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.172"
$ws.Range("E22").Value = "  -1.17%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.003"
$ws.Range("E23").Value = "  -0.77%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "148.15"
$ws.Range("E24").Value = "  -0.55%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1245"
$ws.Range("E25").Value = "  -6.00%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.594"
$ws.Range("E26").Value = "  -3.98%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.96"
$ws.Range("E27").Value = "  +1.23%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06229"
$ws.Range("E28").Value = "  -0.97%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.360"
$ws.Range("E29").Value = "  -1.70%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.279"
$ws.Range("E30").Value = "  -4.01%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.468"
$ws.Range("E31").Value = "  -3.82%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.433"
$ws.Range("E32").Value = "  -4.55%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.621"
$ws.Range("E33").Value = "  -3.73%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9938"
$ws.Range("E34").Value = "  -4.61%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6042"
$ws.Range("E35").Value = "  -1.56%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.403"
$ws.Range("E36").Value = "  -0.46%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.712"
$ws.Range("E37").Value = "  +0.17%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.112"
$ws.Range("E38").Value = "  -1.48%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01613"
$ws.Range("E39").Value = "  -1.73%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8669"
$ws.Range("E40").Value = "  -2.47%  "

$ws.Range("D41").Value = "1.071.24"
$ws.Range("E41").Value = "  -3.99%  "

$ws.Range("E42").Value = "  -1.14%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.91"
$ws.Range("E43").Value = "  -1.94%  "

$ws.Range("D44").Value = "1.821.06"
$ws.Range("E44").Value = "  -1.30%  "

$ws.Range("E45").Value = "  +0.89%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.00"
$ws.Range("E46").Value = "  -2.60%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.000"
$ws.Range("E47").Value = "  -1.38%  "

$ws.Range("E48").Value = "  -0.81%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.967"
$ws.Range("E49").Value = "  -2.94%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4253"
$ws.Range("E50").Value = "  -1.18%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.981"
$ws.Range("E51").Value = "  -1.92%  "

